# Fixed tests for fund_units
# Rewrites the Investors sheet sample data (TSTF* fund/portfolio-company test
# rows) in place of the old AC / Investor 37 / Ego Pvt Ltd sample rows, and
# extends the table down to row 12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header row (unchanged text, kept for completeness / robustness)
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Name *"
$ws.Range("B1").Value = "PAN"
$ws.Range("C1").Value = "Primary Email *"
$ws.Range("D1").Value = "Tags"
$ws.Range("E1").Value = "Category *"
$ws.Range("F1").Value = "City"

# ---------------------------------------------------------------------
# 2. Data rows 2-10 : fund investors (LP)
# ---------------------------------------------------------------------
$funds = @(
    @{ row = 2;  name = "TSTF1";  pan = "TSTFU1111D"; email = "tstf11@emp.com" },
    @{ row = 3;  name = "TSTF2";  pan = "TSTFU1212D"; email = "tstf12@emp.com" },
    @{ row = 4;  name = "TSTF7";  pan = "TSTFU1112D"; email = "tstf13@emp.com" },
    @{ row = 5;  name = "TSTF8";  pan = "TSTFU1113D"; email = "tstf14@emp.com" },
    @{ row = 6;  name = "TSTF9";  pan = "TSTFU1114D"; email = "tstf15@emp.com" },
    @{ row = 7;  name = "TSTF10"; pan = "TSTFU1115D"; email = "tstf16@emp.com" },
    @{ row = 8;  name = "TSTF11"; pan = "TSTFU1116D"; email = "tstf17@emp.com" },
    @{ row = 9;  name = "TSTF12"; pan = "TSTFU1117D"; email = "tstf18@emp.com" },
    @{ row = 10; name = "TSTF13"; pan = "TSTFU1118D"; email = "tstf19@emp.com" }
)

foreach ($f in $funds) {
    $r = $f.row
    $ws.Cells.Item($r, 1).Value = $f.name
    $ws.Cells.Item($r, 2).Value = $f.pan
    $ws.Cells.Item($r, 3).Value = $f.email
    $ws.Cells.Item($r, 5).Value = "LP"
}

# ---------------------------------------------------------------------
# 3. Data rows 11-12 : portfolio companies
# ---------------------------------------------------------------------
$portcos = @(
    @{ row = 11; name = "TSTF1 Port Co 1"; pan = "TSTFP1111D"; email = "tstfportco11@emp.com"; tag = "Fintech, Deal lead 1" },
    @{ row = 12; name = "TSTF1 Port Co 2"; pan = "TSTFP1212D"; email = "tstfportco12@emp.com"; tag = "Fintech, Deal lead 2" }
)

foreach ($p in $portcos) {
    $r = $p.row
    $ws.Cells.Item($r, 1).Value = $p.name
    $ws.Cells.Item($r, 2).Value = $p.pan
    $ws.Cells.Item($r, 3).Value = $p.email
    $ws.Cells.Item($r, 4).Value = $p.tag
    $ws.Cells.Item($r, 5).Value = "Portfolio Company"
}

# ---------------------------------------------------------------------
# 4. Hyperlinks on column C (mailto: links to each row's own email)
# ---------------------------------------------------------------------
$allRows = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12)
foreach ($r in $allRows) {
    $addr = "mailto:" + $ws.Cells.Item($r, 3).Value
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 3), $addr)
}

# ---------------------------------------------------------------------
# 5. Style the "Tags" blanks (D4:D10) like the Hyperlink column (blank,
#    quote-prefixed Hyperlink-derived style in the source workbook)
# ---------------------------------------------------------------------
$ws.Range("D4:D10").Style = "Hyperlink"

# ---------------------------------------------------------------------
# 6. Column widths - re-run best-fit autosize for the data columns, then
#    widen column D (Tags) manually since it's not best-fit in the target.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(5).AutoFit()
$ws.Columns.Item(6).AutoFit()
$ws.Columns.Item(4).ColumnWidth = 19.5

# ---------------------------------------------------------------------
# 7. Sheet view / selection tweaks
# ---------------------------------------------------------------------
$ws.Range("A13").Select()
